$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SKU cell (A2) is cleared out to a blank/space placeholder value
$ws.Range("A2").Value = " "

# Large Tag Quantity (I2) bumped from 0 to 1
$ws.Range("I2").Value = 1

# Move/restore the active selection to A2 (was E2)
[void]$ws.Range("A2").Select()
